# 8.3.1.2.xlsx -- add a "2020" data column (K) to the SDG 8.3.1.2 table.
#
# Target shape (per the OOXML diff):
#   K3 = 2020                (header, same look as the other year headers)
#   K4 = 2.8218550629805335  (Small enterprises row)
#   K5 = 1.3005071159823327  (Medium-sized enterprises row)
# K4/K5 use a brand-new number format (#,##0.0) and a brand-new 9pt
# "Kyrghyz Times" font; K4 also needs a brand-new "medium top border only"
# border while K5 re-uses the workbook's existing "medium bottom border".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K3: header cell, same formatting family as the existing year headers
# (bold 10pt Times New Roman, medium border top+bottom) -----------------
$ws.Range("I3").Copy()
$ws.Range("K3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K3").Font.Bold = $true
$ws.Range("K3").Font.Size = 10
$ws.Range("K3").Font.Name = "Times New Roman"
$ws.Range("K3").VerticalAlignment = -4107   # xlBottom (sheet default)
$ws.Range("K3").Value = 2020

# --- K4: "Small enterprises" 2020 value ---------------------------------
$ws.Range("J5").Copy()
$ws.Range("K4").PasteSpecial(-4122)   # xlPasteFormats (brings in the medium bottom border)
$ws.Range("K4").Font.Name = "Kyrghyz Times"
$ws.Range("K4").Font.Size = 9
$ws.Range("K4").NumberFormat = "#,##0.0"
$ws.Range("K4").HorizontalAlignment = -4152   # xlRight
$ws.Range("K4").VerticalAlignment = -4107     # xlBottom (sheet default)
$ws.Range("K4").Borders.Item(9).LineStyle = 0 # clear the bottom border copied from J5
$ws.Range("K4").Borders.Item(8).Weight = -4138 # xlMedium top border
$ws.Range("K4").Value = 2.8218550629805335

# --- K5: "Medium-sized enterprises" 2020 value --------------------------
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)   # xlPasteFormats (keeps the medium bottom border)
$ws.Range("K5").Font.Name = "Kyrghyz Times"
$ws.Range("K5").Font.Size = 9
$ws.Range("K5").NumberFormat = "#,##0.0"
$ws.Range("K5").HorizontalAlignment = -4152   # xlRight
$ws.Range("K5").VerticalAlignment = -4107     # xlBottom (sheet default)
$ws.Range("K5").Value = 1.3005071159823327

# Move the UI selection cursor to match the saved workbook state.
$ws.Range("L8").Select()
